$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 306 (pushes the existing rows 306..323 down to 307..324)
$ws.Rows.Item(306).Insert()

# Populate the new weekly record in row 306
$ws.Range("A306").Value = 5
$ws.Range("B306").Value = "Macroferia Regional de Talca"
$ws.Range("C306").Value = "Maule"
$ws.Range("D306").Value = 44706
$ws.Range("E306").Value = 7
$ws.Range("F306").Value = 100114014
$ws.Range("G306").Value = "Betarraga"
$ws.Range("H306").Value = "Sin especificar"
$ws.Range("I306").Value = "Primera"
$ws.Range("J306").Value = 5000
$ws.Range("K306").Value = 600
$ws.Range("L306").Value = 600
$ws.Range("M306").Value = 600
$ws.Range("N306").Value = "`$/paquete 5 unidades"
$ws.Range("O306").Value = "Región del Maule"
$ws.Range("P306").Value = 120
$ws.Range("Q306").Value = 5
$ws.Range("R306").Value = "Hortaliza"
